$d = $word.ActiveDocument

# Locate the paragraph "Ver no Jupiter Salvar em pdf Salvar em docx".
# It sits between a blank paragraph (to be removed) and the
# "© 2020 ..." copyright paragraph (to be removed), right after the
# "LOQ4233: Gestão de Negócios (Requisito fraco)" paragraph.

$target = "Ver no Jupiter Salvar em pdf Salvar em docx"

$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $target) {
        $idx = $i
        break
    }
}

if ($idx -gt 0) {
    # Delete from the highest index down to the lowest so earlier
    # deletions don't shift the indices of paragraphs we still need
    # to remove.
    $d.Paragraphs.Item($idx + 1).Range.Delete()  # "© 2020 ..." paragraph
    $d.Paragraphs.Item($idx).Range.Delete()      # "Ver no Jupiter ..." paragraph
    $d.Paragraphs.Item($idx - 1).Range.Delete()  # blank paragraph before it
}
